$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 419 (pushes existing rows 419.. down by 2)
$ws.Rows.Item(419).Insert()
$ws.Rows.Item(419).Insert()

# Row 419: new "Primera" quality record for date 2023-07-17 (serial 45124)
$ws.Cells.Item(419, 1).Value = 1
$ws.Cells.Item(419, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(419, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(419, 4).Value = 45124
$ws.Cells.Item(419, 5).Value = 15
$ws.Cells.Item(419, 6).Value = 100114014
$ws.Cells.Item(419, 7).Value = "Betarraga"
$ws.Cells.Item(419, 8).Value = "Sin especificar"
$ws.Cells.Item(419, 9).Value = "Primera"
$ws.Cells.Item(419, 10).Value = 850
$ws.Cells.Item(419, 11).Value = 700
$ws.Cells.Item(419, 12).Value = 800
$ws.Cells.Item(419, 13).Value = 747
$ws.Cells.Item(419, 14).Value = "$/paquete 4 unidades"
$ws.Cells.Item(419, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(419, 16).Value = 187
$ws.Cells.Item(419, 17).Value = 4
$ws.Cells.Item(419, 18).Value = "Hortaliza"

# Row 420: new "Segunda" quality record for date 2023-07-17 (serial 45124)
$ws.Cells.Item(420, 1).Value = 1
$ws.Cells.Item(420, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(420, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(420, 4).Value = 45124
$ws.Cells.Item(420, 5).Value = 15
$ws.Cells.Item(420, 6).Value = 100114014
$ws.Cells.Item(420, 7).Value = "Betarraga"
$ws.Cells.Item(420, 8).Value = "Sin especificar"
$ws.Cells.Item(420, 9).Value = "Segunda"
$ws.Cells.Item(420, 10).Value = 1050
$ws.Cells.Item(420, 11).Value = 700
$ws.Cells.Item(420, 12).Value = 800
$ws.Cells.Item(420, 13).Value = 743
$ws.Cells.Item(420, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(420, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(420, 16).Value = 149
$ws.Cells.Item(420, 17).Value = 5
$ws.Cells.Item(420, 18).Value = "Hortaliza"
